$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells for columns L (Got) and M (On order) ---
$ws.Range("L1").Value = "got"
$ws.Range("M1").Value = "onorder"

# --- Update J2:J19 into a shared formula (same formula already present, just re-applied) ---
$ws.Range("J2:J19").Formula = "=I2/H2*E2"
# Re-applying the formula to J17 introduces a spurious inherited number format; reset it.
$ws.Range("J17").Style = "Normal"

# --- Stock / order tracking data (columns L, M, N) ---
# Row 14 noted first ("lots" in stock)
$ws.Range("L14").Value = "lots"

# Row 19: part swapped out for a new display, with stock tracking info
$ws.Range("B19").Value = "0.96"" 80x160 IPS TFT"
$ws.Range("L19").Value = "sipeed langan nano"
$ws.Range("M19").Value = 5
$ws.Range("N19").Value = "aliexpress"

# Remaining rows
$ws.Range("L2").Value = 18

$ws.Range("L3").Value = 6

$ws.Range("L4").Value = 3
$ws.Range("M4").Value = 20
$ws.Range("N4").Value = "aliexpress"

$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 10
$ws.Range("N5").Value = "farnell"

$ws.Range("L6").Value = 2
$ws.Range("M6").Value = 10
$ws.Range("N6").Value = "farnell"

$ws.Range("L7").Value = 98

$ws.Range("L8").Value = 15

$ws.Range("L9").Value = 6

$ws.Range("L10").Value = 23

$ws.Range("L11").Value = 2
$ws.Range("M11").Value = 10
$ws.Range("N11").Value = "farnell"

$ws.Range("L12").Value = 13

$ws.Range("L13").Value = 3

$ws.Range("L15").Value = 2
$ws.Range("M15").Value = 6
$ws.Range("N15").Value = "aliexpress"

$ws.Range("L16").Value = 7

$ws.Range("L17").Value = "lots"

$ws.Range("L18").Value = 2
$ws.Range("M18").Value = 10
$ws.Range("N18").Value = "farnell"

# Update the active selection to match where editing left off
[void]$ws.Range("O18").Select()
